# LOB1019.docx edit: insert manual line breaks (<w:br/>) inside several
# long runs of text, splitting them at natural boundaries (sentence /
# list-item boundaries), without altering any of the actual characters.
#
# Word's Find/Replace turns the "^l" replacement code into a manual line
# break (w:br) when inserted into a run, which is exactly the edit the
# diff describes (a single <w:r> whose <w:t> content gets interleaved
# with <w:br/> elements).

$d = $word.ActiveDocument

function Break-Between($before, $after) {
    $find = $before + $after
    $replace = $before + "^l" + $after
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1) English "Programa resumido" summary sentence.
Break-Between "ideal gas," "temperature, heat and the laws of thermodynamics."

# 2) Portuguese "Programa" (detailed syllabus) - split into 8 pieces at the
#    7 numbered-item / sentence boundaries.
Break-Between "tensão superficial, capilaridade;" "2) Dinâmica de fluidos"
Break-Between "lei de Hagen-Poiseuille;" "3) Oscilações"
Break-Between "amortecido e forçado, ressonância; " "4) Ondas"
Break-Between "batimentos, efeito Doppler;" "5) Temperatura e calor"
Break-Between "gás ideal e graus de liberdade;" "6) Termodinâmica"
Break-Between "processos reversíveis e" "irreversíveis, entropia"
Break-Between "máquinas térmicas e" "eficiência."

# 3) English "Programa" (detailed syllabus) - split into 7 pieces at the
#    6 numbered-item boundaries (item 1 is itself split once more).
Break-Between "Archimedes’ principle, surface" "tension and capillarity;"
Break-Between "tension and capillarity;" "2) Fluids in motion"
Break-Between "and the Hagen-Poiseuille law;" "3) Oscillation"
Break-Between "damped and forced oscillations, resonance;" "4) Waves"
Break-Between "intensity and sound level, beats, Doppler effect;" "5) Temperature and heat"
Break-Between "degrees of freedom for an ideal gas;" "6) Thermodynamics"

# 4) Bibliography - split into 5 references.
Break-Between "Edgard Blucher (2008)." "RESNICK, R.; HALLIDAY, D."
Break-Between "Fundamentos de Física. Vol.2, LTC (2008)." "TIPLER, P.; MOSCA, G."
Break-Between "Física para Cientistas e Engenheiros. Vol.2, LTC (2008)." "SEARS, F. W."
Break-Between "Pearson Addison Wesley (2009)." "JEWETT Jr, John W."
